# Applies the "comparison of cross-platform frameworks" edit described in the
# task diff against Literature_review.docx.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Learning-curve bullet: append a follow-up question.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "learning curve: Am I familiar with the programming language? If not how much effort is required to learn the programming language?",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "learning curve: Am I familiar with the programming language? If not how much effort is required to learn the programming language? Can the development team at Coach in a Box easily grasp the language or framework?",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Interpreted apps paragraph: "languages or technologies" -> "programming
#    languages or frameworks".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "independently of the platform using languages or technologies like Ruby",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "independently of the platform using programming languages or frameworks like Ruby",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "when framework supports them" -> "when the framework supports them".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "can be available to apps only when framework supports them (Xanthopoulos",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "can be available to apps only when the framework supports them (Xanthopoulos",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4) "These frameworks are written" -> "Frameworks used to develop these apps
#    are written".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "These frameworks are written in one of several programming languages and compiled",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Frameworks used to develop these apps are written in one of several programming languages and compiled",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Drop "in the case of " before "a suggestion to correct a deprecated
#    call."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "specific needs e.g. in the case of a suggestion to correct a deprecated call.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "specific needs e.g. a suggestion to correct a deprecated call.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Move the (hidden) _GoBack bookmark, and insert the new
#    "2.2.4.3 Comparison of some cross-platform mobile app frameworks"
#    section right after the "(iii) Generated mobile applications" content
#    (i.e. after the paragraph ending "...Xinogalos, 2013, p. 216).").
# ---------------------------------------------------------------------------

# Locate the paragraph that now ends with the Applause/Xanthopoulos citation.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*XText framework*") {
        $targetPara = $candidate
        break
    }
}

$cursor = $targetPara.Range
$cursor.InsertParagraphAfter() | Out-Null

# -- blank paragraph --------------------------------------------------------
$idx = $targetPara.Index + 1
$cursor = $d.Paragraphs.Item($idx).Range

# -- underlined heading paragraph -------------------------------------------
$cursor.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$headingRange = $d.Paragraphs.Item($idx).Range
$headingRange.InsertAfter("2.2.4.3 Comparison of some cross-platform mobile app frameworks")
$headingRange.Font.Underline = 1

# -- criteria paragraph (+ relocated _GoBack bookmark) ----------------------
$headingPara = $d.Paragraphs.Item($idx)
$headingPara.Range.InsertParagraphAfter() | Out-Null
$idx = $idx + 1
$criteriaRange = $d.Paragraphs.Item($idx).Range
$criteriaRange.InsertAfter("Some frameworks used in cross-platform mobile app development will be compared based on the following criteria: price, learning curve, performance, platform, stability and online resource and support.")

# Re-fetch the paragraph/range after inserting text, then drop the bookmark
# at the very end of it (this both creates the new _GoBack and removes it
# from its old location, matching Word's "single _GoBack" behaviour).
$criteriaPara = $d.Paragraphs.Item($idx)
$bmRange = $criteriaPara.Range.Duplicate
$bmRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# -- 14 trailing blank paragraphs --------------------------------------------
$lastPara = $d.Paragraphs.Item($idx)
for ($n = 0; $n -lt 14; $n++) {
    $lastPara.Range.InsertParagraphAfter() | Out-Null
    $idx = $idx + 1
    $lastPara = $d.Paragraphs.Item($idx)
}
